# "Modifications du profil Administrateur"
#
# In the "Gestion du menu administrateur" block (second copy of the task
# table, starting at row 25) two task rows progress along their
# A faire / En cours / Terminé tracking columns (B / C / D):
#
#   - Row 29 "Mettre à jour mon profil" moves from "A faire" (col B) to
#     "Terminé" (col D), and gets a start date + an end date recorded
#     (both 2016-05-24, serial 42514), reusing the same date formatting
#     already used for the other finished rows (e.g. F9/F10/F28).
#   - Row 32 "Rechercher d'autres anciens élèves" moves from "A faire"
#     (col B) to "En cours" (col C).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# --- Row 29 : "Mettre à jour mon profil" -> Terminé, with start/end dates
$ws.Range("B29").ClearContents()
$ws.Range("D29").Value = "X"

# Pull the date number format (incl. border/alignment) from an existing
# finished-task date cell so the new cells are styled identically, then
# write the actual date values (Excel serial date 42514 = 2016-05-24).
$ws.Range("F9").Copy()
$ws.Range("E29:F29").PasteSpecial(-4122)
$ws.Range("E29").Value = 42514
$ws.Range("F29").Value = 42514

# --- Row 32 : "Rechercher d'autres anciens élèves" -> En cours
$ws.Range("B32").ClearContents()
$ws.Range("C32").Value = "X"

# Reflect the author's final on-screen selection/scroll position.
$win = $excel.ActiveWindow
$win.ScrollRow = 22
$win.ScrollColumn = 1
$ws.Range("F30").Select()
